$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(55,5).Style = $ws.Cells.Item(54,5).Style
$ws.Cells.Item(55,6).Style = $ws.Cells.Item(54,6).Style
